# Append the updated list of FTSE 100 tickers (duplicate block, excluding "PSH")
# to the bottom of the worksheet, extending the used range from A1:A130 to A1:A192.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tickers = @("AAF","ABDN","ABF","ANTO","AUTO","AV","BARC","BATS","BDEV","BEZ","BF.B","BKG","BNZL","BRBY","BRK.B","BT-A","CCH","CRDA","DCC","DGE","ENT","EXPN","FCIT","FRAS","FRES","GLEN","HLMA","HSBA","HSX","IMB","INF","ITRK","JMAT","KGF","LGEN","LLOY","LSEG","MNDI","MNG","OCDO","PHNX","PSON","REL","RMV","RR","RS1","SBRY","SDR","SGE","SGRO","SKG","SMDS","SMT","SN","SPX","SSE","STAN","STJ","ULVR","UU","WEIR","WTB")

$startRow = 131
for ($i = 0; $i -lt $tickers.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $tickers[$i]
}
